# Applies the "Updated capital structure database" edit to the Colombia
# Bank (Money Center) dataset: refreshes all metric columns for the existing
# companies (rows 2-8), re-sorts/relabels a few company rows, clears a couple
# of now-stale debt_ebitda/net_debt_ebitda cells, and appends a new company
# (Banco Davivienda) as row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Colombia"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "7"
$ws.Range("C2").Value = "Bank (Money Center)"
$ws.Range("D2").Value = 0.0496
$ws.Range("E2").Value = -0.00958
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.0002682156833781273
$ws.Range("J2").Value = 0.0002051195126585329
$ws.Range("K2").Value = 2028.4
$ws.Range("L2").Value = 0.1214239962646138
$ws.Range("M2").Value = 1488.8
$ws.Range("N2").Value = 0.04548023375520465
$ws.Range("O2").Value = 0.733977519226977
$ws.Range("P2").Value = 1370.1
$ws.Range("Q2").Value = 0.04185415654755904
$ws.Range("R2").Value = 0.6754584894498128
$ws.Range("S2").Value = 118.7
$ws.Range("T2").Value = 0.07972864051585168
$ws.Range("U2").Value = 29633.4
$ws.Range("V2").Value = 0.9052484947350092
$ws.Range("W2").Value = 0.09324104234527687
$ws.Range("X2").Value = 0.08145004043388948
$ws.Range("Y2").Value = 0.01179100191138739
$ws.Range("Z2").Value = 0.305421698250767
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04517285862784731
$ws.Range("AC2").Value = -0.04517285862784731
$ws.Range("AD2").Value = 59324.4
$ws.Range("AE2").Value = 5.297150938000226
$ws.Range("AF2").Value = 59329.697150938
$ws.Range("AG2").Value = 29696.297150938
$ws.Range("AH2").Value = 0.6444341266908827
$ws.Range("AI2").Value = 0.6471093055792877
$ws.Range("AJ2").Value = 0.475662863657246
$ws.Range("AK2").Value = 0.478580429493952
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 10708.37545126354
$ws.Range("AP2").Value = 5360.342446017689

# Row 3
$ws.Range("A3").Value = "Colombia"
$ws.Range("B3").Value = "Banco Comercial AV Villas S.A. (BVC:VILLAS)"
$ws.Range("C3").Value = "Bank (Money Center)"
$ws.Range("D3").Value = 0.0638
$ws.Range("E3").Value = 0.049
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 45.8
$ws.Range("L3").Value = 0.18490109002826
$ws.Range("M3").Value = 18.4
$ws.Range("N3").Value = 0.05199208816049732
$ws.Range("O3").Value = 0.4017467248908297
$ws.Range("P3").Value = 18.4
$ws.Range("Q3").Value = 0.05199208816049732
$ws.Range("R3").Value = 0.4017467248908297
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 230
$ws.Range("V3").Value = 0.6499011020062165
$ws.Range("W3").Value = 0.09324104234527687
$ws.Range("X3").Value = 0.04659247447585589
$ws.Range("Y3").Value = 0.04664856786942098
$ws.Range("Z3").Value = 0.3608682983682984
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04342856582273767
$ws.Range("AC3").Value = -0.04342856582273767
$ws.Range("AD3").Value = 96.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 96.5
$ws.Range("AG3").Value = -133.5
$ws.Range("AH3").Value = 0.2142539964476022
$ws.Range("AI3").Value = 0.1705248277080756
$ws.Range("AJ3").Value = -0.6057168784029039
$ws.Range("AK3").Value = -0.397439714200655
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0

# Row 4
$ws.Range("A4").Value = "Colombia"
$ws.Range("B4").Value = "Banco Bilbao Vizcaya Argentaria Colombia S.A. (BVC:BBVACOL)"
$ws.Range("C4").Value = "Bank (Money Center)"
$ws.Range("D4").Value = 0.309
$ws.Range("E4").Value = -0.00958
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 144.2
$ws.Range("L4").Value = 0.03591442305297501
$ws.Range("M4").Value = 89.90000000000001
$ws.Range("N4").Value = 0.0699284380833852
$ws.Range("O4").Value = 0.6234396671289876
$ws.Range("P4").Value = 89.90000000000001
$ws.Range("Q4").Value = 0.0699284380833852
$ws.Range("R4").Value = 0.6234396671289876
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 1199.5
$ws.Range("V4").Value = 0.9330273802115744
$ws.Range("W4").Value = 0.09984075330609983
$ws.Range("X4").Value = 0.08905011569306433
$ws.Range("Y4").Value = 0.0107906376130355
$ws.Range("Z4").Value = 1.089904720540731
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04431303937517304
$ws.Range("AC4").Value = -0.04431303937517304
$ws.Range("AD4").Value = 2902.1
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 2902.1
$ws.Range("AG4").Value = 1702.6
$ws.Range("AH4").Value = 0.6930057071901043
$ws.Range("AI4").Value = 0.6777440448388603
$ws.Range("AJ4").Value = 0.5697744461548758
$ws.Range("AK4").Value = 0.5523438767234388
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

# Row 5
$ws.Range("A5").Value = "Colombia"
$ws.Range("B5").Value = "Grupo Aval Acciones y Valores S.A. (BVC:GRUPOAVAL)"
$ws.Range("C5").Value = "Bank (Money Center)"
$ws.Range("D5").Value = 0.06619999999999999
$ws.Range("E5").Value = 0.0495
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 630.9
$ws.Range("L5").Value = 0.1439228031754722
$ws.Range("M5").Value = 404.6
$ws.Range("N5").Value = 0.05326627873298402
$ws.Range("O5").Value = 0.6413060706926613
$ws.Range("P5").Value = 345.8
$ws.Range("Q5").Value = 0.04552515864030122
$ws.Range("R5").Value = 0.5481058804881915
$ws.Range("S5").Value = 58.80000000000001
$ws.Range("T5").Value = 0.1453287197231834
$ws.Range("U5").Value = 10130.7
$ws.Range("V5").Value = 1.333723900050028
$ws.Range("W5").Value = 0.1135917610413928
$ws.Range("X5").Value = 0.08975256762432032
$ws.Range("Y5").Value = 0.02383919341707251
$ws.Range("Z5").Value = 0.3934973653737399
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.04432875065529163
$ws.Range("AC5").Value = -0.04432875065529163
$ws.Range("AD5").Value = 17396.1
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 17396.1
$ws.Range("AG5").Value = 7265.399999999998
$ws.Range("AH5").Value = 0.6960695265265946
$ws.Range("AI5").Value = 0.6541042966238395
$ws.Range("AJ5").Value = 0.4888838048071488
$ws.Range("AK5").Value = 0.4412740060493421
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0

# Row 6
$ws.Range("A6").Value = "Colombia"
$ws.Range("B6").Value = "Grupo Bolívar S.A. (BVC:GRUBOLIVAR)"
$ws.Range("C6").Value = "Bank (Money Center)"
$ws.Range("D6").Value = 0.035
$ws.Range("E6").Value = -0.0348
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 129.1
$ws.Range("L6").Value = 0.08773956775859724
$ws.Range("M6").Value = 29.1
$ws.Range("N6").Value = 0.01882398602755676
$ws.Range("O6").Value = 0.2254066615027111
$ws.Range("P6").Value = 29.1
$ws.Range("Q6").Value = 0.01882398602755676
$ws.Range("R6").Value = 0.2254066615027111
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 2841.2
$ws.Range("V6").Value = 1.837893783556504
$ws.Range("W6").Value = 0.04479683542107637
$ws.Range("X6").Value = 0.1729118180159206
$ws.Range("Y6").Value = -0.1281149825948443
$ws.Range("Z6").Value = 0.1659242887267561
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.04517285862784731
$ws.Range("AC6").Value = -0.04517285862784731
$ws.Range("AD6").Value = 9549.9
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 9549.9
$ws.Range("AG6").Value = 6708.7
$ws.Range("AH6").Value = 0.860677012923809
$ws.Range("AI6").Value = 0.6922976548624451
$ws.Range("AJ6").Value = 0.8127226031546047
$ws.Range("AK6").Value = 0.612482083025207
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0

# Row 7
$ws.Range("A7").Value = "Colombia"
$ws.Range("B7").Value = "Banco de Bogotá S.A. (BVC:BOGOTA)"
$ws.Range("C7").Value = "Bank (Money Center)"
$ws.Range("D7").Value = 0.0496
$ws.Range("E7").Value = 0.102
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.001574228730377329
$ws.Range("J7").Value = 0.001212921558325409
$ws.Range("K7").Value = 618
$ws.Range("L7").Value = 0.2171316140819338
$ws.Range("M7").Value = 439.5
$ws.Range("N7").Value = 0.06009763301472699
$ws.Range("O7").Value = 0.7111650485436893
$ws.Range("P7").Value = 379.6
$ws.Range("Q7").Value = 0.0519068520873501
$ws.Range("R7").Value = 0.6142394822006473
$ws.Range("S7").Value = 59.89999999999998
$ws.Range("T7").Value = 0.1362912400455062
$ws.Range("U7").Value = 7628.4
$ws.Range("V7").Value = 1.043114411125241
$ws.Range("W7").Value = 0.1069703840894535
$ws.Range("X7").Value = 0.06699225171761321
$ws.Range("Y7").Value = 0.03997813237184034
$ws.Range("Z7").Value = 0.4754525089944834
$ws.Range("AA7").Value = 0.0005766865981193144
$ws.Range("AB7").Value = 0.04851925474682969
$ws.Range("AC7").Value = -0.04794256814871038
$ws.Range("AD7").Value = 8962.6
$ws.Range("AE7").Value = 5.297150938000226
$ws.Range("AF7").Value = 8967.897150938001
$ws.Range("AG7").Value = 1339.497150938001
$ws.Range("AH7").Value = 0.550819895599658
$ws.Range("AI7").Value = 0.5976407886210829
$ws.Range("AJ7").Value = 0.1548086808586472
$ws.Range("AK7").Value = 0.1815750997352236
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = 0
$ws.Range("AN7").Value = 1617.797833935018
$ws.Range("AP7").Value = 241.7864893389894

# Row 8
$ws.Range("A8").Value = "Colombia"
$ws.Range("B8").Value = "Bancolombia S.A. (BVC:BCOLOMBIA)"
$ws.Range("C8").Value = "Bank (Money Center)"
$ws.Range("D8").Value = -0.0018
$ws.Range("E8").Value = -0.162
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 262.7
$ws.Range("L8").Value = 0.1072726530278901
$ws.Range("M8").Value = 399
$ws.Range("N8").Value = 0.04054836841089013
$ws.Range("O8").Value = 1.51884278644842
$ws.Range("P8").Value = 399
$ws.Range("Q8").Value = 0.04054836841089013
$ws.Range("R8").Value = 1.51884278644842
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 4824.1
$ws.Range("V8").Value = 0.4902490828345241
$ws.Range("W8").Value = 0.03384698636843868
$ws.Range("X8").Value = 0.06529435556570881
$ws.Range("Y8").Value = -0.03144736919727013
$ws.Range("Z8").Value = 0.1596144069453679
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0.04828531210506485
$ws.Range("AC8").Value = -0.04828531210506485
$ws.Range("AD8").Value = 11285.7
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 11285.7
$ws.Range("AG8").Value = 6461.6
$ws.Range("AH8").Value = 0.5342140889339102
$ws.Range("AI8").Value = 0.596161769832969
$ws.Range("AJ8").Value = 0.3963758381027745
$ws.Range("AK8").Value = 0.4580583418991245
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0

# Row 9
$ws.Range("A9").Value = "Colombia"
$ws.Range("B9").Value = "Banco Davivienda S.A. (BVC:PFDAVVNDA)"
$ws.Range("C9").Value = "Bank (Money Center)"
$ws.Range("D9").Value = 0.0472
$ws.Range("E9").Value = -0.05019999999999999
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 197.7
$ws.Range("L9").Value = 0.1529948924315121
$ws.Range("M9").Value = 108.3
$ws.Range("N9").Value = 0.02255921011519154
$ws.Range("O9").Value = 0.5477996965098635
$ws.Range("P9").Value = 108.3
$ws.Range("Q9").Value = 0.02255921011519154
$ws.Range("R9").Value = 0.5477996965098635
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 2779.5
$ws.Range("V9").Value = 0.5789780656987523
$ws.Range("W9").Value = 0.0560135996600085
$ws.Range("X9").Value = 0.08145004043388948
$ws.Range("Y9").Value = -0.02543644077388098
$ws.Range("Z9").Value = 0.1437694704049844
$ws.Range("AA9").Value = 0
$ws.Range("AB9").Value = 0.04911203098428361
$ws.Range("AC9").Value = -0.04911203098428361
$ws.Range("AD9").Value = 9131.5
$ws.Range("AE9").Value = 0
$ws.Range("AF9").Value = 9131.5
$ws.Range("AG9").Value = 6352
$ws.Range("AH9").Value = 0.6554241254073297
$ws.Range("AI9").Value = 0.7299127126230976
$ws.Range("AJ9").Value = 0.5695481811579258
$ws.Range("AK9").Value = 0.6527659312088296
$ws.Range("AL9").Value = 0
$ws.Range("AM9").Value = 0

# Clear cells that no longer have data (debt_ebitda / net_debt_ebitda were
# dropped for these rows in the refreshed dataset)
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()
$ws.Range("AN8").ClearContents()
$ws.Range("AP8").ClearContents()

